$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are numeric-looking strings (e.g. "63.40", "209.52").
# The source workbook stores every Price/Volume cell as literal text (to keep exact
# formatting such as trailing zeros), so we briefly force text interpretation via
# NumberFormat while writing those specific cells, then restore the default style.
$textForcedCells = @("D5","D8","D11","D14","D15","D17","D18","D20","D23","D24","D25","D28","D43","D45","D46","D48")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.756.28"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.594.63"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "209.52"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  -1.85%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "22.36"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "0.0869"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").Value = "1.820.94"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").Value = "1.599.88"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "3.84"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "0.533"
$ws.Range("E15").Value = "  -3.09%  "
$ws.Range("D16").Value = "27.754.26"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "63.40"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "219.10"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0696"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "7.37"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").Value = "9.75"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").Value = "153.71"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  +4.09%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "15.17"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("D33").Value = "1.381.81"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "64.56"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("E44").Value = "  +2.70%  "
$ws.Range("D45").Value = "1.76"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "5.24"
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").Value = "1.732.07"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("D48").Value = "86.01"
$ws.Range("E48").Value = "  -3.38%  "
$ws.Range("E49").Value = "  +4.66%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("E51").Value = "  -1.07%  "

foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
